# Auto-generated edit script: updates cryptocurrency price/volume data
# scraped from coinranking.com (cryptos list refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number (e.g. "241.23").
# These must be written as text (matching the original inline-string
# cells) instead of letting Excel coerce them into Doubles, which
# would silently change the stored precision/formatting.
$textAddrs = @('D5', 'D6', 'D8', 'D12', 'D14', 'D17', 'D19', 'D20', 'D21', 'D22', 'D23', 'D24', 'D25', 'D26', 'D27', 'D28', 'D29', 'D30', 'D32', 'D34', 'D35', 'D36', 'D37', 'D38', 'D39', 'D40', 'D41', 'D43', 'D44', 'D46', 'D47', 'D48', 'D49', 'D50', 'D51')
foreach ($addr in $textAddrs) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '96.300.86'
$ws.Range("E2").Value = '  +4.53%  '
$ws.Range("D3").Value = '3.672.96'
$ws.Range("E3").Value = '  +10.06%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = '241.23'
$ws.Range("E5").Value = '  +4.28%  '
$ws.Range("D6").Value = '645.48'
$ws.Range("E6").Value = '  +4.78%  '
$ws.Range("E7").Value = '  +5.90%  '
$ws.Range("D8").Value = '0.403'
$ws.Range("E8").Value = '  +4.64%  '
$ws.Range("E9").Value = '  -0.22%  '
$ws.Range("E10").Value = '  +5.69%  '
$ws.Range("D11").Value = '3.671.27'
$ws.Range("E11").Value = '  +10.10%  '
$ws.Range("D12").Value = '43.75'
$ws.Range("E12").Value = '  +1.58%  '
$ws.Range("E13").Value = '  +3.38%  '
$ws.Range("D14").Value = '6.37'
$ws.Range("E14").Value = '  +3.62%  '
$ws.Range("D15").Value = '4.362.69'
$ws.Range("E15").Value = '  +9.99%  '
$ws.Range("D16").Value = '96.183.74'
$ws.Range("E16").Value = '  +4.63%  '
$ws.Range("D17").Value = '0.0000256'
$ws.Range("E17").Value = '  +5.30%  '
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '3.663.42'
$ws.Range("E18").Value = '  +9.81%  '
$ws.Range("B19").Value = 'Uniswap'
$ws.Range("C19").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D19").Value = '13.49'
$ws.Range("E19").Value = '  +23.50%  '
$ws.Range("D20").Value = '8.03'
$ws.Range("E20").Value = '  -0.86%  '
$ws.Range("D21").Value = '18.46'
$ws.Range("E21").Value = '  +5.80%  '
$ws.Range("D22").Value = '519.47'
$ws.Range("E22").Value = '  +4.99%  '
$ws.Range("D23").Value = '0.490'
$ws.Range("E23").Value = '  +11.63%  '
$ws.Range("D24").Value = '3.45'
$ws.Range("E24").Value = '  +0.38%  '
$ws.Range("D25").Value = '0.0000198'
$ws.Range("E25").Value = '  +8.20%  '
$ws.Range("D26").Value = '6.75'
$ws.Range("E26").Value = '  +6.19%  '
$ws.Range("D27").Value = '97.57'
$ws.Range("E27").Value = '  +8.54%  '
$ws.Range("D28").Value = '12.60'
$ws.Range("E28").Value = '  +5.85%  '
$ws.Range("D29").Value = '3.17'
$ws.Range("E29").Value = '  +21.51%  '
$ws.Range("D30").Value = '11.66'
$ws.Range("E30").Value = '  +4.48%  '
$ws.Range("E31").Value = '  +2.64%  '
$ws.Range("D32").Value = '1.00'
$ws.Range("E32").Value = '  -0.08%  '
$ws.Range("E33").Value = '  +4.79%  '
$ws.Range("D34").Value = '0.996'
$ws.Range("E34").Value = '  +0.00%  '
$ws.Range("D35").Value = '31.95'
$ws.Range("E35").Value = '  +12.72%  '
$ws.Range("D36").Value = '0.580'
$ws.Range("E36").Value = '  +9.43%  '
$ws.Range("D37").Value = '568.30'
$ws.Range("E37").Value = '  +0.45%  '
$ws.Range("D38").Value = '7.87'
$ws.Range("E38").Value = '  +6.76%  '
$ws.Range("D39").Value = '1.49'
$ws.Range("E39").Value = '  +8.71%  '
$ws.Range("D40").Value = '0.951'
$ws.Range("E40").Value = '  +9.27%  '
$ws.Range("D41").Value = '0.153'
$ws.Range("E41").Value = '  +2.95%  '
$ws.Range("E42").Value = '  -0.09%  '
$ws.Range("B43").Value = 'Filecoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D43").Value = '5.78'
$ws.Range("E43").Value = '  +6.92%  '
$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").Value = '0.0430'
$ws.Range("E44").Value = '  +3.78%  '
$ws.Range("E45").Value = '  +2.82%  '
$ws.Range("D46").Value = '23.77'
$ws.Range("E46").Value = '  +0.31%  '
$ws.Range("D47").Value = '2.23'
$ws.Range("E47").Value = '  +5.62%  '
$ws.Range("D48").Value = '54.19'
$ws.Range("E48").Value = '  +4.97%  '
$ws.Range("D49").Value = '3.51'
$ws.Range("E49").Value = '  -2.74%  '
$ws.Range("D50").Value = '8.23'
$ws.Range("E50").Value = '  +3.29%  '
$ws.Range("D51").Value = '3.14'
$ws.Range("E51").Value = '  +3.81%  '

# Restore the default (unstyled) cell style now that the text is
# safely stored, so no stray number-format style lingers on the cell.
foreach ($addr in $textAddrs) {
    $ws.Range($addr).Style = "Normal"
}
